$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.523.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.943.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.72%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.40'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.73%  '

$ws.Range("E9").Value = '  -1.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '55.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.06%  '

$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.820'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.229.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.66%  '

$ws.Range("E17").Value = '  -2.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.947.43'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.433.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0863'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("E25").Value = '  -3.85%  '

$ws.Range("E26").Value = '  +0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.131'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.13%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.66%  '

$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("E32").Value = '  +1.87%  '

$ws.Range("E33").Value = '  -2.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0627'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.48%  '

$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("E38").Value = '  -2.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.69%  '

$ws.Range("E40").Value = '  -0.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0981'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("E42").Value = '  +0.94%  '

$ws.Range("E43").Value = '  -3.27%  '

$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.348.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.30%  '

$ws.Range("E47").Value = '  -4.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.45%  '
